$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 value updates ---
$ws.Range("A2").Value = 2177602
$ws.Range("B2").Value = 11931

# --- Clear rows 3-6: drop A:D entirely, keep E cells but strip value+style ---
$ws.Range("A3:D6").Clear()
$ws.Range("E3:E7").ClearContents()
$ws.Range("E3:E7").Style = "Normal"

# --- Conditional formatting: shrink range to A2 only, bump priorities by 10 ---
$fc = $ws.Range("A2:A6").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($ws.Range("A2"))
}
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).Priority = 30 + $i
}

# --- Selection ---
$ws.Range("D6").Select()
